# Apply the 3.3.2 data-update edit: add column R (year 2021) to the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column Q (year 2020) onto the new column R
# so every row in R3:R33 picks up the same per-row style as its Q neighbor.
$ws.Range("Q3:Q33").Copy()
$ws.Range("R3:R33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the 2021 values.
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 58.14349653559799
$ws.Range("R5").Value = 50.405857641278807
$ws.Range("R6").Value = 65.995789757646122
$ws.Range("R7").Value = 47.339416388110941
$ws.Range("R8").Value = 44.18457369250482
$ws.Range("R9").Value = 50.379263611270765
$ws.Range("R10").Value = 54.819947539591084
$ws.Range("R11").Value = 47.679920417302263
$ws.Range("R12").Value = 61.861274529713718
$ws.Range("R13").Value = 36.712395096811576
$ws.Range("R14").Value = 26.872053459579295
$ws.Range("R15").Value = 46.638444428499682
$ws.Range("R16").Value = 51.155081745820631
$ws.Range("R17").Value = 43.08338023862634
$ws.Range("R18").Value = 58.934228062068456
$ws.Range("R19").Value = 54.51979816984521
$ws.Range("R20").Value = 52.474443936678909
$ws.Range("R21").Value = 56.519551395440942
$ws.Range("R22").Value = 46.970408642555192
$ws.Range("R23").Value = 27.43769048802011
$ws.Range("R24").Value = 66.104415920267911
$ws.Range("R25").Value = 88.246666265390886
$ws.Range("R26").Value = 71.914698721605745
$ws.Range("R27").Value = 105.10059183863845
$ws.Range("R28").Value = 63.980940123966526
$ws.Range("R29").Value = 55.546587096180644
$ws.Range("R30").Value = 73.505198287622903
$ws.Range("R31").Value = 43.916363725083563
$ws.Range("R32").Value = 40.980198843051781
$ws.Range("R33").Value = 47.015458682814909

# Leave the selection where the author left it after the edit.
$ws.Range("T3").Select()
